$wb = $excel.ActiveWorkbook

# The "optimization_parameters" sheet had a stray leftover row (row 16:
# "Sheet", 3, 4) that is no longer needed - remove it. This shifts the
# "simulation_timepoints" row (old row 17) up to row 16, and drops the
# now-unused "Sheet" shared string / number-format style along with it.
$wsParams = $wb.Worksheets.Item("optimization_parameters")
$wsParams.Rows.Item(16).Select()
$wsParams.Rows.Item(16).Delete()

# Move the active tab from "optimization_parameters" to "threshold_b".
$wsThreshold = $wb.Worksheets.Item("threshold_b")
$wsThreshold.Activate()
$wsThreshold.Range("A2").Select()
